# [Phatnttse][Update Export To Excel]
# Translate the Supplier List export template header row (row 1) to
# Vietnamese and fix the District/Ward column order in the placeholder
# row (row 2) so it lines up with the (already correct) header order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: human-readable column headers -> Vietnamese -------------------
$ws.Range("A1").Value = "Danh mục"
$ws.Range("B1").Value = "Tên nhà cung cấp"
$ws.Range("C1").Value = "Đường dẫn"
$ws.Range("D1").Value = "Số điện thoại"
$ws.Range("E1").Value = "Mô tả"
$ws.Range("F1").Value = "Giá"
$ws.Range("G1").Value = "Địa chỉ chi tiết"
$ws.Range("H1").Value = "Tỉnh/Thành phố"
$ws.Range("I1").Value = "Quận/Huyện"
$ws.Range("J1").Value = "Phường/Xã"
$ws.Range("K1").Value = "Website"
$ws.Range("L1").Value = "Thời gian phản hồi"
$ws.Range("M1").Value = "Chỉ số ưu tiên"
$ws.Range("N1").Value = "Nổi bật"
$ws.Range("O1").Value = "Thời gian kết thúc ưu tiên"
$ws.Range("P1").Value = "Trạng thái duyệt"
$ws.Range("Q1").Value = "Trạng thái hoàn thành"
$ws.Range("R1").Value = "Giảm giá"
$ws.Range("S1").Value = "Đánh giá"

# --- Row 2: merge-field placeholders ---------------------------------------
# Columns I and J swap so District/Ward line up with the new header order.
$ws.Range("I2").Value = "[[%Field:District%]]"
$ws.Range("J2").Value = "[[%Field:Ward%]]"

# --- Misc view state --------------------------------------------------------
$ws.Range("R19").Select() | Out-Null
